$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "302.97"
Set-TextValue $ws.Range("E2") "1.48%"
Set-TextValue $ws.Range("D3") "32.13"
Set-TextValue $ws.Range("E3") "2.88%"
Set-TextValue $ws.Range("D4") "4.928"
Set-TextValue $ws.Range("E4") "-3.74%"
Set-TextValue $ws.Range("D5") "0.07833"
Set-TextValue $ws.Range("E5") "-1.43%"
Set-TextValue $ws.Range("D6") "2.014"
Set-TextValue $ws.Range("E6") "-9.69%"
Set-TextValue $ws.Range("D7") "7.831"
Set-TextValue $ws.Range("E7") "0.70%"
Set-TextValue $ws.Range("D8") "3.819"
Set-TextValue $ws.Range("E8") "-1.17%"
Set-TextValue $ws.Range("D9") "0.9208"
Set-TextValue $ws.Range("E9") "-0.14%"
Set-TextValue $ws.Range("D10") "0.1758"
Set-TextValue $ws.Range("E10") "1.80%"
Set-TextValue $ws.Range("D11") "0.07901"
Set-TextValue $ws.Range("E11") "6.07%"
Set-TextValue $ws.Range("D12") "0.08588"
Set-TextValue $ws.Range("E12") "-7.02%"
Set-TextValue $ws.Range("D13") "0.03162"
Set-TextValue $ws.Range("E13") "3.17%"
Set-TextValue $ws.Range("D14") "0.1004"
Set-TextValue $ws.Range("E14") "0.13%"
Set-TextValue $ws.Range("D15") "0.001520"
Set-TextValue $ws.Range("E15") "0.44%"
Set-TextValue $ws.Range("D16") "0.005832"
Set-TextValue $ws.Range("E16") "-3.22%"
Set-TextValue $ws.Range("E17") "2,108.87%"
Set-TextValue $ws.Range("D18") "3.464"
Set-TextValue $ws.Range("E18") "-0.44%"
Set-TextValue $ws.Range("D19") "2.158"
Set-TextValue $ws.Range("E19") "-4.86%"
Set-TextValue $ws.Range("D22") "4.275"
Set-TextValue $ws.Range("E22") "9.26%"
Set-TextValue $ws.Range("E23") "17.16%"
Set-TextValue $ws.Range("D24") "0.04568"
Set-TextValue $ws.Range("E24") "-1.12%"
Set-TextValue $ws.Range("D25") "0.001224"
Set-TextValue $ws.Range("E25") "-1.78%"
Set-TextValue $ws.Range("D26") "0.004450"
Set-TextValue $ws.Range("E26") "-0.67%"
Set-TextValue $ws.Range("E27") "4.22%"
Set-TextValue $ws.Range("D39") "0.01743"
Set-TextValue $ws.Range("E39") "-0.74%"
Set-TextValue $ws.Range("D40") "0.04780"
Set-TextValue $ws.Range("E40") "3.91%"
Set-TextValue $ws.Range("D41") "0.007589"
Set-TextValue $ws.Range("E41") "8.87%"
Set-TextValue $ws.Range("D42") "0.1366"
Set-TextValue $ws.Range("E42") "0.32%"
Set-TextValue $ws.Range("D43") "0.002360"
Set-TextValue $ws.Range("E43") "7.82%"
Set-TextValue $ws.Range("D44") "0.01058"
Set-TextValue $ws.Range("E44") "5.08%"
Set-TextValue $ws.Range("D45") "0.00006324"
Set-TextValue $ws.Range("E45") "0.16%"
Set-TextValue $ws.Range("E46") "0.06%"
Set-TextValue $ws.Range("E47") "-61.11%"
Set-TextValue $ws.Range("D48") "0.8205"
Set-TextValue $ws.Range("E48") "9.93%"
Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("E49") "0.06%"
Set-TextValue $ws.Range("D50") "0.0002000"
Set-TextValue $ws.Range("E50") "0.06%"

Write-Host "Applied crypto price/volume updates"
